$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 1
$tr.Text = "Scrum ceremonies let us coordinate between ourselves and where we needed to be in the project. Post-scrum time was a reliable time for us all to work and talk together at the same time, instead of asynchronously. "

# Paragraph 2
[void]$tr.InsertAfter([char]13)
[void]$tr.InsertAfter("Material results included the Sprint documentation in Sprint1.md. ")

# Paragraph 3 (built from several runs)
[void]$tr.InsertAfter([char]13)
[void]$tr.InsertAfter("The most valuable results were process results, which included concrete listing of tasks in ")
[void]$tr.InsertAfter("ZenHub")
[void]$tr.InsertAfter(" for us to complete and orientation")
[void]$tr.InsertAfter("/planning ")
[void]$tr.InsertAfter("to complete those tasks. ")
